# Append two new rows (4 and 5) to the "Shahbaz Nadeem" sheet, duplicating
# the existing Dubai (row 3) and Abu Dhabi (row 2) match entries.
# All values in this sheet are stored as text (numberStoredAsText), so every
# cell here is forced to text via NumberFormat "@" before assignment and then
# restored to the default "Normal" style so no new per-cell style index is
# introduced (matches the rest of the sheet, which uses the default style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 4: duplicate of the "Dubai (DSC)" match (currently row 3)
Set-TextValue "A4" " Dubai (DSC)"
Set-TextValue "B4" " October 13 2020"
Set-TextValue "C4" "Super Kings won by 20 runs"
Set-TextValue "D4" "Sunrisers Hyderabad"
Set-TextValue "E4" "Chennai Super Kings"
Set-TextValue "F4" ("Shahbaz Nadeem" + $nbsp)
Set-TextValue "G4" "5"
Set-TextValue "H4" "5"
Set-TextValue "I4" "1"
Set-TextValue "J4" "0"
Set-TextValue "K4" "100.00"

# Row 5: duplicate of the "Abu Dhabi" match (currently row 2)
Set-TextValue "A5" " Abu Dhabi"
Set-TextValue "B5" " November 08 2020"
Set-TextValue "C5" "Capitals won by 17 runs"
Set-TextValue "D5" "Sunrisers Hyderabad"
Set-TextValue "E5" "Delhi Capitals"
Set-TextValue "F5" ("Shahbaz Nadeem" + $nbsp)
Set-TextValue "G5" "2"
Set-TextValue "H5" "3"
Set-TextValue "I5" "0"
Set-TextValue "J5" "0"
Set-TextValue "K5" "66.66"
